# Reverts the "Lastoppdata (#13)" change:
#  1. For rows 39-51, columns A (IndID) and B (Register) had been swapped;
#     swap them back so column A holds the per-row indicator id (nra1..nra13)
#     and column B holds the constant register name "nra".
#  2. Restore the sheet view: selection should be the full rows 4:7
#     (activeCell A4, sqref A4:XFD7), with the view scrolled back to the
#     top-left (A1), instead of being scrolled to A39 with B39 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 39; $r -le 51; $r++) {
    $colA = $ws.Cells.Item($r, 1).Value2
    $colB = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value = $colB
    $ws.Cells.Item($r, 2).Value = $colA
}

# Reset the window scroll position to the top-left (A1), then select whole
# rows 4:7, matching activeCell="A4" sqref="A4:XFD7" with no topLeftCell
# override.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$null = $ws.Rows("4:7").Select()
